# Fill header for review
# Reorders the Kategorie/Merkmal value pairs in rows 12, 15, 17, 22 and 23
# of the "Dataset" sheet so that the (category, feature) pairs appear in a
# new column order, while keeping the same set of values per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Row 12
$ws.Range("C12").Value = "Strassen_und_Gehsteige"
$ws.Range("D12").Value = "GehsteigbreiteMin"
$ws.Range("E12").Value = "Lage_Gelaende_Planzeichen"
$ws.Range("F12").Value = "AnFluchtlinie"

# Row 15
$ws.Range("C15").Value = "Flaeche"
$ws.Range("D15").Value = "Flaechen"
$ws.Range("E15").Value = "Dach"
$ws.Range("F15").Value = "DachneigungMax"

# Row 17
$ws.Range("C17").Value = "Flaeche"
$ws.Range("D17").Value = "Flaechen"
$ws.Range("E17").Value = "Dach"
$ws.Range("F17").Value = "Dachart"
$ws.Range("G17").Value = "Dach"
$ws.Range("H17").Value = "BegruenungDach"
$ws.Range("I17").Value = "Ausgestaltung_und_Sonstiges"
$ws.Range("J17").Value = "GebaeudeBautyp"

# Row 22
$ws.Range("C22").Value = "Ausgestaltung_und_Sonstiges"
$ws.Range("D22").Value = "UnterbrechungGeschlosseneBauweise"
$ws.Range("E22").Value = "Lage_Gelaende_Planzeichen"
$ws.Range("F22").Value = "Planzeichen"

# Row 23
$ws.Range("C23").Value = "Flaeche"
$ws.Range("D23").Value = "Flaechen"
$ws.Range("E23").Value = "Dach"
$ws.Range("F23").Value = "Dachart"
$ws.Range("G23").Value = "Dach"
$ws.Range("H23").Value = "BegruenungDach"
$ws.Range("I23").Value = "Lage_Gelaende_Planzeichen"
$ws.Range("J23").Value = "Planzeichen"
